$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Column C ("Update Date") held dates typed as plain text
#    ("15 05 2025", "22 03 2025", ...). Replace the text with real
#    Excel date serial values and format the whole data column as a
#    date (this also drops the now-unused date strings from the
#    shared string table automatically).
# ------------------------------------------------------------------
$ws.Range("C3").Value = 45792   # 15 05 2025
$ws.Range("C4").Value = 45738   # 22 03 2025
$ws.Range("C5").Value = 45782   # 05 05 2025
$ws.Range("C6").Value = 45747   # 31 03 2025
$ws.Range("C7").Value = 45786   # 09 05 2025
$ws.Range("C8").Value = 45782   # 05 05 2025
$ws.Range("C9").Value = 45790   # 13 05 2025
$ws.Range("C17").Value = 45792  # 15 05 2025

# Apply a real date number format (maps to the built-in date format)
# across the whole column of data rows, including the still-empty ones.
$ws.Range("C3:C21").NumberFormat = "mm-dd-yy"

# ------------------------------------------------------------------
# 2. Columns H and I (and the blank rows of C) had been force-
#    formatted as Text (to stop Excel mangling the typed dates).
#    That is no longer needed, so the forced text format is cleared
#    back to General across the table body.
# ------------------------------------------------------------------
$ws.Range("H3:I21").NumberFormat = "General"

# ------------------------------------------------------------------
# 3. Restore the active selection to H5 (matches the saved view).
# ------------------------------------------------------------------
$ws.Range("H5").Select()
